# Atualização de bases das ligas, do dia: 28-05-2024 às 20:56
#
# The data rows below got re-ordered (the match results/odds for each row were
# shuffled between a few rows while the sequential "id"/"Date" columns in A/D
# stayed put). Re-create that by reading each affected row's B:AD cell values
# first (so the reads all see the original/"before" data), then writing them
# back out according to the new row mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

function Read-RowValues($row) {
    $vals = @()
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals += $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Write-RowValues($row, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, $firstCol + $i).Value2 = $vals[$i]
    }
}

# Snapshot every row that participates in a re-shuffle before writing any of
# them back, since several groups rotate data among more than two rows.
$rowsToRead = @(22, 23, 44, 45, 46, 47, 61, 62, 66, 67)
$orig = @{}
foreach ($r in $rowsToRead) {
    $orig[$r] = Read-RowValues $r
}

# destination row -> source row (which row's original B:AD data now lands here)
$mapping = @{
    22 = 23
    23 = 22
    44 = 45
    45 = 47
    46 = 44
    47 = 46
    61 = 62
    62 = 61
    66 = 67
    67 = 66
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    Write-RowValues $dst $orig[$src]
}
